# 1) Duplicate the existing "2022-Q2" sheet so the OLD Q2 data/styles survive
#    untouched on its own tab, then repurpose the original tab object for the
#    new "2022-Q3" data (this is what lets the new sheet land on sheetId=2
#    and the relocated "2022-Q2" sheet land on sheetId=3, matching the diff).
$wb = $excel.ActiveWorkbook
$wsQ2 = $wb.Worksheets.Item(2)
$wsQ2.Copy($null, $wsQ2)

$wsCopy = $wb.Worksheets.Item(3)
$wsCopy.Name = "2022-Q2_tmp"
$wsQ2.Name = "2022-Q3"
$wsCopy.Name = "2022-Q2"

# 2) Overwrite the (renamed) "2022-Q3" tab with the new quarter's fund table.
$wsQ3 = $wb.Worksheets.Item("2022-Q3")

$wsQ3.Cells.Item(1,1).Value = ""
$wsQ3.Range("B1").Value = "基金代码"
$wsQ3.Range("C1").Value = "基金名称"
$wsQ3.Range("D1").Value = "基金规模"
$wsQ3.Range("E1").Value = "股票总仓位"
$wsQ3.Range("F1").Value = "仓位占比"
$wsQ3.Range("G1").Value = "持有市值(亿元)"
$wsQ3.Range("H1").Value = "仓位排名"

$wsQ3.Range("A2").Value = 0
$wsQ3.Range("B2").Value = "161611"
$wsQ3.Range("C2").Value = "融通内需驱动混合A/B"
$wsQ3.Range("D2").Value = "8.72"
$wsQ3.Range("E2").Value = "90.68"
$wsQ3.Range("F2").Value = "3.99"
$wsQ3.Range("G2").Value = "0.3479"
$wsQ3.Range("H2").Value = 9

$wsQ3.Range("A3").Value = 1
$wsQ3.Range("B3").Value = "014109"
$wsQ3.Range("C3").Value = "融通内需驱动混合C"
$wsQ3.Range("D3").Value = "4.06"
$wsQ3.Range("E3").Value = "90.68"
$wsQ3.Range("F3").Value = "3.99"
$wsQ3.Range("G3").Value = "0.1620"
$wsQ3.Range("H3").Value = 9

$wsQ3.Range("A4").Value = 2
$wsQ3.Range("B4").Value = "001319"
$wsQ3.Range("C4").Value = "农银汇理信息传媒主题股票"
$wsQ3.Range("D4").Value = "2.62"
$wsQ3.Range("E4").Value = "80.55"
$wsQ3.Range("F4").Value = "3.55"
$wsQ3.Range("G4").Value = "0.0930"
$wsQ3.Range("H4").Value = 9

# 3) Insert the new Q3 summary row into "总计", pushing the old Q2 row down.
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Range("A3").EntireRow.Insert()

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q2"
$wsTotal.Range("C3").Value = 2
$wsTotal.Range("D3").Value = 0.05

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 3
$wsTotal.Range("D2").Value = 0.6
